$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dhnfgbf"
$ws.Range("D13").Value = "uklyun"
$ws.Range("F7").Value = "rthsf"
$ws.Range("I10").Value = "rtnrgtsf"

$ws.Range("I10").Select()
